$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 208.17857
$ws.Range("I33").Value = 203.34782
$ws.Range("J33").Value = 230.4
$ws.Range("K33").Value = 203.34782
$ws.Range("L33").Value = 230.4
$ws.Range("M33").Value = 25.65217999999999
$ws.Range("N33").Value = -688.4
$ws.Range("H112").Value = 34483984
$ws.Range("J112").Value = 37038310
$ws.Range("L112").Value = 111114930
$ws.Range("N112").Value = -111117146
$ws.Range("H131").Value = 1539.5883
$ws.Range("I131").Value = 1591.7142
$ws.Range("J131").Value = 1296.3334
$ws.Range("K131").Value = 4775.142599999999
$ws.Range("L131").Value = 3889.0002
$ws.Range("M131").Value = 264.8574000000008
$ws.Range("N131").Value = -13969.0002
$ws.Range("H132").Value = 229760.34
$ws.Range("I132").Value = 288578.44
$ws.Range("J132").Value = 1023.3333
$ws.Range("K132").Value = 865735.3200000001
$ws.Range("L132").Value = 3069.9999
$ws.Range("M132").Value = -863205.3200000001
$ws.Range("N132").Value = -8129.9999
$ws.Range("H135").Value = 2027.9667
$ws.Range("I135").Value = 658.4583
$ws.Range("J135").Value = 7506
$ws.Range("K135").Value = 5926.1247
$ws.Range("L135").Value = 67554
$ws.Range("M135").Value = -3391.1247
$ws.Range("N135").Value = -72624
$ws.Range("H137").Value = 25642556
$ws.Range("I137").Value = 1001.28
$ws.Range("J137").Value = 71431050
$ws.Range("K137").Value = 3003.84
$ws.Range("L137").Value = 214293150
$ws.Range("M137").Value = -453.8400000000001
$ws.Range("N137").Value = -214298250

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4232.2534
$ws.Range("I32").Value = 4579.068
$ws.Range("K32").Value = 4579.068
$ws.Range("M32").Value = -4292.068
$ws.Range("H61").Value = 2861.2
$ws.Range("I61").Value = 1685.3334
$ws.Range("K61").Value = 1685.3334
$ws.Range("M61").Value = -1473.3334
$ws.Range("H132").Value = 2576.5
$ws.Range("I132").Value = 2771.077
$ws.Range("J132").Value = 2346.5454
$ws.Range("K132").Value = 8313.231
$ws.Range("L132").Value = 7039.6362
$ws.Range("M132").Value = -5783.231
$ws.Range("N132").Value = -12099.6362
$ws.Range("H135").Value = 42085.8
$ws.Range("J135").Value = 42085.8
$ws.Range("L135").Value = 42085.8
$ws.Range("N135").Value = -52225.8
$ws.Range("H136").Value = 2861.2
$ws.Range("I136").Value = 1685.3334
$ws.Range("K136").Value = 5056.0002
$ws.Range("M136").Value = -2506.0002
$ws.Range("H139").Value = 52730.715
$ws.Range("J139").Value = 52730.715
$ws.Range("L139").Value = 52730.715
$ws.Range("N139").Value = -63010.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2901.5789
$ws.Range("I31").Value = 1736
$ws.Range("J31").Value = 3317.8572
$ws.Range("K31").Value = 1736
$ws.Range("L31").Value = 3317.8572
$ws.Range("M31").Value = -1441
$ws.Range("N31").Value = -3907.8572
$ws.Range("H34").Value = 2901.5789
$ws.Range("I34").Value = 1736
$ws.Range("J34").Value = 3317.8572
$ws.Range("K34").Value = 1736
$ws.Range("L34").Value = 3317.8572
$ws.Range("M34").Value = -1534
$ws.Range("N34").Value = -3721.8572
$ws.Range("H107").Value = 326.02274
$ws.Range("I107").Value = 337.91177
$ws.Range("J107").Value = 285.6
$ws.Range("K107").Value = 337.91177
$ws.Range("L107").Value = 285.6
$ws.Range("M107").Value = 1582.08823
$ws.Range("N107").Value = -4125.6
$ws.Range("H122").Value = 6580358.5
$ws.Range("I122").Value = 10417856
$ws.Range("K122").Value = 31253568
$ws.Range("M122").Value = -31251118
$ws.Range("H132").Value = 3247.0833
$ws.Range("I132").Value = 2729.0557
$ws.Range("J132").Value = 4801.1665
$ws.Range("K132").Value = 8187.1671
$ws.Range("L132").Value = 14403.4995
$ws.Range("M132").Value = -5657.1671
$ws.Range("N132").Value = -19463.4995
$ws.Range("H134").Value = 1781.9565
$ws.Range("I134").Value = 1926.8334
$ws.Range("J134").Value = 1260.4
$ws.Range("K134").Value = 5780.5002
$ws.Range("L134").Value = 3781.2
$ws.Range("M134").Value = -3245.5002
$ws.Range("N134").Value = -8851.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1272.5883
$ws.Range("J5").Value = 966.6667
$ws.Range("L5").Value = 2900.0001
$ws.Range("N5").Value = -3124.0001
$ws.Range("H122").Value = 20835066
$ws.Range("I122").Value = 47619610
$ws.Range("K122").Value = 428576490
$ws.Range("M122").Value = -428574040
$ws.Range("H132").Value = 45456972
$ws.Range("I132").Value = 100001896
$ws.Range("K132").Value = 900017064
$ws.Range("M132").Value = -900014534
$ws.Range("H133").Value = 6013.593
$ws.Range("I133").Value = 2869.8333
$ws.Range("J133").Value = 6911.8096
$ws.Range("K133").Value = 8609.499899999999
$ws.Range("L133").Value = 20735.4288
$ws.Range("M133").Value = -3549.499899999999
$ws.Range("N133").Value = -30855.4288
$ws.Range("H135").Value = 1272.5883
$ws.Range("J135").Value = 966.6667
$ws.Range("L135").Value = 8700.0003
$ws.Range("N135").Value = -13770.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1577.8334
$ws.Range("I102").Value = 1159.3334
$ws.Range("K102").Value = 1159.3334
$ws.Range("M102").Value = 462.6666
$ws.Range("H122").Value = 4306.852
$ws.Range("I122").Value = 4722.091
$ws.Range("J122").Value = 2479.8
$ws.Range("K122").Value = 14166.273
$ws.Range("L122").Value = 7439.400000000001
$ws.Range("M122").Value = -11716.273
$ws.Range("N122").Value = -12339.4
$ws.Range("H132").Value = 2151.5
$ws.Range("I132").Value = 2064.4375
$ws.Range("J132").Value = 2238.5625
$ws.Range("K132").Value = 6193.3125
$ws.Range("L132").Value = 6715.6875
$ws.Range("M132").Value = -3663.3125
$ws.Range("N132").Value = -11775.6875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5079.6978
$ws.Range("I122").Value = 5847.933
$ws.Range("J122").Value = 3306.8462
$ws.Range("K122").Value = 17543.799
$ws.Range("L122").Value = 9920.5386
$ws.Range("M122").Value = -15093.799
$ws.Range("N122").Value = -14820.5386
$ws.Range("H132").Value = 78340.07000000001
$ws.Range("I132").Value = 114850.8
$ws.Range("J132").Value = 5318.6
$ws.Range("K132").Value = 344552.4
$ws.Range("L132").Value = 15955.8
$ws.Range("N132").Value = -21015.8
$ws.Range("M132").Value = -342022.4
